# Auto-generated Excel COM-interop script
# Applies numeric value updates to cells in the Leve profit tables
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2005.8182
$ws.Range("I62").Value = 1243.5714
$ws.Range("K62").Value = 1243.5714
$ws.Range("M62").Value = -619.5714
$ws.Range("H65").Value = 2005.8182
$ws.Range("I65").Value = 1243.5714
$ws.Range("K65").Value = 6217.857
$ws.Range("M65").Value = -3097.857
$ws.Range("H76").Value = 2698.4
$ws.Range("I76").Value = 2679.6875
$ws.Range("J76").Value = 2773.25
$ws.Range("K76").Value = 2679.6875
$ws.Range("L76").Value = 2773.25
$ws.Range("M76").Value = -2364.6875
$ws.Range("N76").Value = -3403.25
$ws.Range("H79").Value = 2698.4
$ws.Range("I79").Value = 2679.6875
$ws.Range("J79").Value = 2773.25
$ws.Range("K79").Value = 2679.6875
$ws.Range("L79").Value = 2773.25
$ws.Range("M79").Value = -1587.6875
$ws.Range("N79").Value = -4957.25
$ws.Range("H115").Value = 1118.1
$ws.Range("I115").Value = 341.57144
$ws.Range("J115").Value = 2930
$ws.Range("K115").Value = 1024.71432
$ws.Range("L115").Value = 8790
$ws.Range("M115").Value = 542.28568
$ws.Range("N115").Value = -11924
$ws.Range("H131").Value = 4724.533
$ws.Range("I131").Value = 7775
$ws.Range("J131").Value = 3615.2727
$ws.Range("K131").Value = 23325
$ws.Range("L131").Value = 10845.8181
$ws.Range("M131").Value = -18285
$ws.Range("N131").Value = -20925.8181
$ws.Range("H137").Value = 5524.5
$ws.Range("I137").Value = 6049.1113
$ws.Range("J137").Value = 4344.125
$ws.Range("K137").Value = 18147.3339
$ws.Range("L137").Value = 13032.375
$ws.Range("M137").Value = -15597.3339
$ws.Range("N137").Value = -18132.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7814255
$ws.Range("I2").Value = 22729120
$ws.Range("J2").Value = 1707.1428
$ws.Range("K2").Value = 22729120
$ws.Range("L2").Value = 1707.1428
$ws.Range("M2").Value = -22729007
$ws.Range("N2").Value = -1933.1428
$ws.Range("H32").Value = 5303.63
$ws.Range("I32").Value = 5216.638
$ws.Range("K32").Value = 5216.638
$ws.Range("M32").Value = -4929.638
$ws.Range("H61").Value = 2830.606
$ws.Range("I61").Value = 2133.4546
$ws.Range("J61").Value = 3179.182
$ws.Range("K61").Value = 2133.4546
$ws.Range("L61").Value = 3179.182
$ws.Range("M61").Value = -1921.4546
$ws.Range("N61").Value = -3603.182
$ws.Range("H74").Value = 3129.0833
$ws.Range("I74").Value = 2220.7144
$ws.Range("J74").Value = 4400.8
$ws.Range("K74").Value = 2220.7144
$ws.Range("L74").Value = 4400.8
$ws.Range("M74").Value = -1346.7144
$ws.Range("N74").Value = -6148.8
$ws.Range("H77").Value = 3129.0833
$ws.Range("I77").Value = 2220.7144
$ws.Range("J77").Value = 4400.8
$ws.Range("K77").Value = 11103.572
$ws.Range("L77").Value = 22004
$ws.Range("M77").Value = -6735.572
$ws.Range("N77").Value = -30740
$ws.Range("H116").Value = 7814255
$ws.Range("I116").Value = 22729120
$ws.Range("J116").Value = 1707.1428
$ws.Range("K116").Value = 22729120
$ws.Range("L116").Value = 1707.1428
$ws.Range("M116").Value = -22726826
$ws.Range("N116").Value = -6295.1428
$ws.Range("H132").Value = 2109.3704
$ws.Range("I132").Value = 1955.7916
$ws.Range("J132").Value = 3338
$ws.Range("K132").Value = 5867.3748
$ws.Range("L132").Value = 10014
$ws.Range("M132").Value = -3337.3748
$ws.Range("N132").Value = -15074
$ws.Range("H135").Value = 28945.277
$ws.Range("J135").Value = 28945.277
$ws.Range("L135").Value = 28945.277
$ws.Range("N135").Value = -39085.277
$ws.Range("H136").Value = 2830.606
$ws.Range("I136").Value = 2133.4546
$ws.Range("J136").Value = 3179.182
$ws.Range("K136").Value = 6400.3638
$ws.Range("L136").Value = 9537.545999999998
$ws.Range("M136").Value = -3850.3638
$ws.Range("N136").Value = -14637.546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7814255
$ws.Range("I3").Value = 22729120
$ws.Range("J3").Value = 1707.1428
$ws.Range("K3").Value = 22729120
$ws.Range("L3").Value = 1707.1428
$ws.Range("M3").Value = -22729006
$ws.Range("N3").Value = -1935.1428
$ws.Range("H86").Value = 53386
$ws.Range("I86").Value = 5460
$ws.Range("J86").Value = 73925.71000000001
$ws.Range("K86").Value = 5460
$ws.Range("L86").Value = 73925.71000000001
$ws.Range("M86").Value = -4337
$ws.Range("N86").Value = -76171.71000000001
$ws.Range("H89").Value = 53386
$ws.Range("I89").Value = 5460
$ws.Range("J89").Value = 73925.71000000001
$ws.Range("K89").Value = 27300
$ws.Range("L89").Value = 369628.55
$ws.Range("M89").Value = -21684
$ws.Range("N89").Value = -380860.55
$ws.Range("H94").Value = 543.0454999999999
$ws.Range("I94").Value = 459.1875
$ws.Range("J94").Value = 766.6667
$ws.Range("K94").Value = 459.1875
$ws.Range("L94").Value = 766.6667
$ws.Range("M94").Value = -8.1875
$ws.Range("N94").Value = -1668.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2896.3655
$ws.Range("I31").Value = 1992.1892
$ws.Range("K31").Value = 1992.1892
$ws.Range("M31").Value = -1697.1892
$ws.Range("H34").Value = 2896.3655
$ws.Range("I34").Value = 1992.1892
$ws.Range("K34").Value = 1992.1892
$ws.Range("M34").Value = -1790.1892

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5000.5386
$ws.Range("I56").Value = 5000.5386
$ws.Range("K56").Value = 5000.5386
$ws.Range("M56").Value = -4470.5386
$ws.Range("H82").Value = 2666.6667
$ws.Range("H85").Value = 2666.6667
$ws.Range("H88").Value = 4042.1538
$ws.Range("J88").Value = 4042.1538
$ws.Range("L88").Value = 12126.4614
$ws.Range("N88").Value = -12982.4614
$ws.Range("H91").Value = 4042.1538
$ws.Range("J91").Value = 4042.1538
$ws.Range("L91").Value = 12126.4614
$ws.Range("N91").Value = -15090.4614
$ws.Range("H136").Value = 1748.5555
$ws.Range("I136").Value = 1467.125
$ws.Range("K136").Value = 4401.375
$ws.Range("M136").Value = 698.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 30670
$ws.Range("J21").Value = 30670
$ws.Range("L21").Value = 30670
$ws.Range("N21").Value = -31016
$ws.Range("H30").Value = 30670
$ws.Range("J30").Value = 30670
$ws.Range("L30").Value = 30670
$ws.Range("N30").Value = -30880
$ws.Range("H102").Value = 56252.79
$ws.Range("I102").Value = 2575
$ws.Range("J102").Value = 95291.17999999999
$ws.Range("K102").Value = 2575
$ws.Range("L102").Value = 95291.17999999999
$ws.Range("M102").Value = -953
$ws.Range("N102").Value = -98535.17999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 52632490
$ws.Range("I16").Value = 71429576
$ws.Range("J16").Value = 627.4
$ws.Range("K16").Value = 71429576
$ws.Range("L16").Value = 627.4
$ws.Range("M16").Value = -71429406
$ws.Range("N16").Value = -967.4
$ws.Range("H82").Value = 2128.1428
$ws.Range("I82").Value = 1575
$ws.Range("J82").Value = 2983
$ws.Range("K82").Value = 1575
$ws.Range("L82").Value = 2983
$ws.Range("M82").Value = -1214
$ws.Range("N82").Value = -3705
$ws.Range("H85").Value = 2128.1428
$ws.Range("I85").Value = 1575
$ws.Range("J85").Value = 2983
$ws.Range("K85").Value = 1575
$ws.Range("L85").Value = 2983
$ws.Range("M85").Value = -327
$ws.Range("N85").Value = -5479
$ws.Range("H97").Value = 20609.375
$ws.Range("J97").Value = 20609.375
$ws.Range("L97").Value = 20609.375
$ws.Range("N97").Value = -22591.375
$ws.Range("H122").Value = 3030.5483
$ws.Range("I122").Value = 2517.4707
$ws.Range("J122").Value = 3653.5715
$ws.Range("K122").Value = 7552.4121
$ws.Range("L122").Value = 10960.7145
$ws.Range("M122").Value = -5102.4121
$ws.Range("N122").Value = -15860.7145
$ws.Range("H133").Value = 32326
$ws.Range("J133").Value = 32326
$ws.Range("L133").Value = 32326
$ws.Range("N133").Value = -37386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 48277.863
$ws.Range("I126").Value = 64459.875
$ws.Range("J126").Value = 5125.8335
$ws.Range("K126").Value = 193379.625
$ws.Range("L126").Value = 15377.5005
$ws.Range("M126").Value = -190909.625
$ws.Range("N126").Value = -20317.5005
